# ---------------------------------------------------------------------------
# Applies the BIL122HW3 homework-answer edit:
#   1. "2 / (" splits into two runs: "2 / " and "("
#   2. "N^3" -> "3^N" (split into three runs: "...yaklaşık 3", "^N", " adet...budur.")
#   3. The "_GoBack" bookmark moves from right after "eklenebilirdi" to right
#      before the "(" that precedes "7.015.254.043.203.144.209"; the
#      ") " + "civarı" + ". " runs merge into a single ") civarı. " run.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: split the run "2 / (" into "2 / " and "(" without altering text.
# A formatting toggle (set then reset) forces Word to break the run at the
# boundary while leaving the visible text untouched.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("2 / (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $splitPoint = $r1.Start + 4   # right after "2 / ", before "("
    $tail = $d.Range($splitPoint, $r1.End)
    $tail.Bold = 1
    $tail.Bold = 0
}

# ---------------------------------------------------------------------------
# Change 2: "yaklaşık N^3 adet" -> "yaklaşık 3^N adet", split across 3 runs.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("N^3", $true, $false, $false, $false, $false, $true, 1, $false, "3^N", 2)

if ($found2) {
    $r2b = $d.Content
    $needle = "yaklaşık 3^N adet fonksiyon çağırılıyor. Yinelemeli olanda ise her adımda 4 işlem yapılıyor, bu da taş başına 4*N adet işlem demek. Diğerine kıyasla ne kadar küçük kaldığı bariz. Süre farkının sebebi de budur."
    $found2b = $r2b.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2b) {
        $runStart = $r2b.Start
        $runEnd = $r2b.End

        $split1 = $runStart + 10   # right after "yaklaşık 3"
        $split2 = $runStart + 12   # right after "^N"

        $part2 = $d.Range($split1, $runEnd)
        $part2.Bold = 1
        $part2.Bold = 0

        $part3 = $d.Range($split2, $runEnd)
        $part3.Bold = 1
        $part3.Bold = 0
    }
}

# ---------------------------------------------------------------------------
# Change 3a: merge ") " + "civarı" + ". " into a single ") civarı. " run,
# while protecting the neighboring "7.015...209" and "Bundan sonraki..." runs
# from being swept into the same merge (Word/engine auto-merges
# identically-formatted adjacent runs whenever an edit touches a run
# boundary, so we give the neighbors a transient, then-reverted, formatting
# difference to shield them).
# ---------------------------------------------------------------------------
$rNumProtect = $d.Content
$foundNumProtect = $rNumProtect.Find.Execute("7.015.254.043.203.144.209", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundNumProtect) {
    $rNumProtect.Bold = 1
}

$rAfterProtect = $d.Content
$foundAfterProtect = $rAfterProtect.Find.Execute("Bundan sonraki adımın 9 kentilyonu geçeceği bariz. Bu durumda veri kaybı yaşanır.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundAfterProtect) {
    $rAfterProtect.Bold = 1
}

$rMerge = $d.Content
$foundMerge = $rMerge.Find.Execute(") civarı. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundMerge) {
    $rMerge.Text = "CIVARI_PLACEHOLDER"
}
$rMerge2 = $d.Content
$foundMerge2 = $rMerge2.Find.Execute("CIVARI_PLACEHOLDER", $true, $false, $false, $false, $false, $true, 1, $false, ") civarı. ", 2)

$rNumUnprotect = $d.Content
$foundNumUnprotect = $rNumUnprotect.Find.Execute("7.015.254.043.203.144.209", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundNumUnprotect) {
    $rNumUnprotect.Bold = 0
}

$rAfterUnprotect = $d.Content
$foundAfterUnprotect = $rAfterUnprotect.Find.Execute("Bundan sonraki adımın 9 kentilyonu geçeceği bariz. Bu durumda veri kaybı yaşanır.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundAfterUnprotect) {
    $rAfterUnprotect.Bold = 0
}

# ---------------------------------------------------------------------------
# Change 3b: move the "_GoBack" bookmark from right after "eklenebilirdi" to
# right before the "(" preceding "7.015.254.043.203.144.209".
# ---------------------------------------------------------------------------
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}

$rBm = $d.Content
$foundBm = $rBm.Find.Execute("adımda hesaplanan taş sayısı 7 kentilyon (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundBm) {
    $bmPoint = $rBm.End - 1   # right before "("
    $bmRange = $d.Range($bmPoint, $bmPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Output "done"
